$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# The "Status" text is a shared label that shows up on three sheets; every
# cell currently reading "Ready for handoff" flips to the handback message.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ===================== zh-cn =====================
$wsZh = $wb.Worksheets.Item("zh-cn")

# Existing hyperlink addresses on this sheet, captured in sheet order:
# [0]=A2 (source md), [1]=B2 (.md), [2]=D2 (handoff xlf),
# [3]=A3 (source md), [4]=B3 (.md), [5]=D3 (handoff xlf)
$zhAddrs = @()
foreach ($h in $wsZh.Hyperlinks) {
    $zhAddrs += $h.Address
}

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhAddrs[0], $null, $null, $wsZh.Range("A2").Text) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhAddrs[2], $null, $null, $wsZh.Range("D2").Text) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhAddrs[3], $null, $null, $wsZh.Range("A3").Text) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhAddrs[5], $null, $null, $wsZh.Range("D3").Text) | Out-Null

$wsZh.Range("F2").Font.Underline = $true
$wsZh.Range("F2").Font.Color = 15570276
$wsZh.Range("G2").Font.Underline = $true
$wsZh.Range("G2").Font.Color = 15570276
$wsZh.Range("F3").Font.Underline = $true
$wsZh.Range("F3").Font.Color = 15570276
$wsZh.Range("G3").Font.Underline = $true
$wsZh.Range("G3").Font.Color = 15570276

$wsZh.Range("H2").Value = "2016-03-18 17:27:21"
$wsZh.Range("H3").Value = "2016-03-18 17:27:21"

# ===================== de-de =====================
$wsDe = $wb.Worksheets.Item("de-de")

$deAddrs = @()
foreach ($h in $wsDe.Hyperlinks) {
    $deAddrs += $h.Address
}

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deAddrs[0], $null, $null, $wsDe.Range("A2").Text) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deAddrs[2], $null, $null, $wsDe.Range("D2").Text) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deAddrs[3], $null, $null, $wsDe.Range("A3").Text) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deAddrs[5], $null, $null, $wsDe.Range("D3").Text) | Out-Null

$wsDe.Range("F2").Font.Underline = $true
$wsDe.Range("F2").Font.Color = 15570276
$wsDe.Range("G2").Font.Underline = $true
$wsDe.Range("G2").Font.Color = 15570276
$wsDe.Range("F3").Font.Underline = $true
$wsDe.Range("F3").Font.Color = 15570276
$wsDe.Range("G3").Font.Underline = $true
$wsDe.Range("G3").Font.Color = 15570276

$wsDe.Range("H2").Value = "2016-03-18 17:27:35"
$wsDe.Range("H3").Value = "2016-03-18 17:27:35"
